$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (Food): amount 750 -> 900, date -> 2025-08-15
$ws.Range("B2").Value = 900
$ws.Range("C2").Value = 45884.22928240741

# Propagate C2's date formatting (style) down to the new date cells C3:C5
# by copying the cell (copies format + value), then overwrite each value.
$ws.Range("C2").Copy($ws.Range("C3:C5"))

# New row 3: Meds
$ws.Range("A3").Value = "Meds"
$ws.Range("B3").Value = 550
$ws.Range("C3").Value = 45884.22928240741

# New row 4: Travel
$ws.Range("A4").Value = "Travel"
$ws.Range("B4").Value = 750
$ws.Range("C4").Value = 45883.22928240741

# New row 5: Misc.
$ws.Range("A5").Value = "Misc."
$ws.Range("B5").Value = 70
$ws.Range("C5").Value = 45883.22928240741
